$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-recorded event row (row 14) ---
$ws.Range("A14").Value = 0.6645833333333333
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 57

# --- (Re)apply the elapsed-time formula down through the new row ---
$ws.Range("D8:D14").FormulaR1C1 = "=RC[-2]*60+RC[-1]+R5C4"

# --- Move the active selection to the cell that was last edited ---
$ws.Range("F14").Select()
